$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values that were stored in columns C and D for every data row (rows 1-10).
$swap = @{
    1  = @(1, 5)
    2  = @(2, 20)
    3  = @(3, 30)
    4  = @(4, 20)
    5  = @(5, 20)
    6  = @(6, 30)
    7  = @(7, 5)
    8  = @(8, 30)
    9  = @(9, 30)
    10 = @(10, 5)
}

foreach ($row in $swap.Keys) {
    $vals = $swap[$row]
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("D$row").Value = $vals[1]
}

# Update formula in E1 (standalone formula, no longer references column A).
$ws.Range("E1").Formula = '=CONCATENATE("insert into [dbo].[payment] values(",B1,",",C1,",",D1,")")'

# Update the shared formula for E2:E10 (also drops the column A reference).
$ws.Range("E2:E10").Formula = '=CONCATENATE("insert into [dbo].[payment] values(",B2,",",C2,",",D2,")")'

# Move the active selection from G14 to E14.
$ws.Range("E14").Select()
